# Update sheet (tab) names and stimulus file names to reflect the new
# timestamped filenames used for more accurate stimulus presentation
# time-logging.

$wb = $excel.ActiveWorkbook

# --- Rename worksheets ---------------------------------------------------
$wb.Worksheets.Item(1).Name = "GNG_TO-16512555473261945"
$wb.Worksheets.Item(2).Name = "NB_TO-16512555494852378"
$wb.Worksheets.Item(3).Name = "RS_TO-16512555494872408"
$wb.Worksheets.Item(4).Name = "TOL_TO-1651255549533237"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16512555496102386"

# --- Sheet 1: GNG ----------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-1651255547295194.csv"
$ws1.Range("B3").Value = "GNG_stims-16512555473081934.csv"
$ws1.Range("B4").Value = "go_stims-1651255547310194.csv"
$ws1.Range("B5").Value = "GNG_stims-1651255547325203.csv"

# --- Sheet 2: NB -------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16512555477472394.csv"
$ws2.Range("B3").Value = "ZB-match_6-16512555476032395.csv"
$ws2.Range("B4").Value = "ZB-match_5-16512555476952374.csv"
$ws2.Range("B5").Value = "TB-16512555494632428.csv"
$ws2.Range("B6").Value = "TB-1651255549429237.csv"
$ws2.Range("B7").Value = "ZB-match_5-16512555475042367.csv"
$ws2.Range("B8").Value = "OB-16512555487702377.csv"
$ws2.Range("B9").Value = "OB-1651255548212238.csv"
$ws2.Range("B10").Value = "TB-16512555491082363.csv"

# --- Sheet 3: RS (no content change, only rename handled above) --------

# --- Sheet 4: TOL --------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-1651255549500239.csv"
$ws4.Range("B3").Value = "ZM_stims-1651255549489237.csv"
$ws4.Range("B4").Value = "MM_stims-1651255549516238.csv"
$ws4.Range("B5").Value = "ZM_stims-1651255549501238.csv"
$ws4.Range("B6").Value = "MM_stims-16512555495322378.csv"
$ws4.Range("B7").Value = "ZM_stims-1651255549517237.csv"

# --- Sheet 5: vSAT -------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16512555495632381.csv"
$ws5.Range("B3").Value = "SAT_stims-16512555495382383.csv"
$ws5.Range("B4").Value = "vSAT_stims-16512555495782382.csv"
$ws5.Range("B5").Value = "vSAT_stims-1651255549595238.csv"
